$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.942.57"
$ws.Range("E2").Value = "  +10.94%  "
$ws.Range("D3").Value = "3.262.74"
$ws.Range("E3").Value = "  +6.36%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'397.26"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").Value = "'109.85"
$ws.Range("E6").Value = "  +8.00%  "
$ws.Range("D7").Value = "'0.561"
$ws.Range("E7").Value = "  +5.22%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "  +6.79%  "
$ws.Range("D10").Value = "'39.26"
$ws.Range("E10").Value = "  +6.14%  "
$ws.Range("D11").Value = "'0.0970"
$ws.Range("E11").Value = "  +14.32%  "
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("D13").Value = "3.769.85"
$ws.Range("E13").Value = "  +5.87%  "
$ws.Range("D14").Value = "'8.11"
$ws.Range("E14").Value = "  +5.51%  "
$ws.Range("E15").Value = "  +3.97%  "
$ws.Range("D16").Value = "3.257.82"
$ws.Range("E16").Value = "  +5.88%  "
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").Value = "'10.82"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").Value = "56.768.66"
$ws.Range("E19").Value = "  +10.54%  "
$ws.Range("E20").Value = "  +4.83%  "
$ws.Range("E21").Value = "  +10.91%  "
$ws.Range("D22").Value = "'12.88"
$ws.Range("E22").Value = "  +5.22%  "
$ws.Range("D23").Value = "'306.36"
$ws.Range("E23").Value = "  +15.87%  "
$ws.Range("D24").Value = "'75.11"
$ws.Range("E24").Value = "  +7.45%  "
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").Value = "'28.31"
$ws.Range("E26").Value = "  +5.30%  "
$ws.Range("D27").Value = "'7.94"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  +4.90%  "
$ws.Range("D29").Value = "'7.26"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("E30").Value = "  +3.44%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").Value = "'37.40"
$ws.Range("E34").Value = "  +3.91%  "
$ws.Range("D35").Value = "'0.0481"
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("E36").Value = "  +3.54%  "
$ws.Range("D37").Value = "'51.58"
$ws.Range("E37").Value = "  +3.01%  "
$ws.Range("E38").Value = "  +25.01%  "
$ws.Range("D39").Value = "'3.56"
$ws.Range("E39").Value = "  +7.60%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").Value = "'134.60"
$ws.Range("E41").Value = "  +5.35%  "
$ws.Range("D42").Value = "'1.93"
$ws.Range("E42").Value = "  +4.86%  "
$ws.Range("D43").Value = "'17.35"
$ws.Range("E43").Value = "  +4.42%  "
$ws.Range("E44").Value = "  +4.44%  "
$ws.Range("D45").Value = "'3.94"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("D47").Value = "'22.01"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").Value = "2.152.73"
$ws.Range("E48").Value = "  +4.09%  "
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("D50").Value = "'2.05"
$ws.Range("E50").Value = "  +43.24%  "
$ws.Range("D51").Value = "'2.38"
$ws.Range("E51").Value = "  -4.43%  "
